$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# --- Replace the full Tasks paragraphs FIRST, so the short substring replacements
#     below ("SSL", "flask") do not collide with occurrences inside these texts. ---
ReplaceText "Because SSL is a startup, Wout helped with building the ground on which they are still building. He made the template of how a flask micro service will look in the future and integrated them in a Kubernetes environment." "Migration of data from exam simulator from old to new Database transforming, Flask, cleaning and updating the data to match the new Data Models implemented for the new application version. Technologies: Python, Git, SQL, Microsoft Azure SQL Databases"

ReplaceText "In this school project Wout made a website in Elixir with the Phoenix framework. He had already learned the basics of web, but in this project it was the first time he worked with a sophisticated framework." "Data Anonymatizator App to encrypth and anonymaize confidential from SQL Databases from FeenPOP. Technologies: Python, Git, Streamlit , Threading, Mathematics, SQL, Pandas"

# --- Simple field replacements (first project block: SSL internship) ---
ReplaceText " SSL" " SOLITA INTERNAL"
ReplaceText "09/2020 - 02/2021" "5/2022 - 5/2022"
ReplaceText "INTERNSHIP" "DABASE MIGRATION - EXAM SIMULATOR"
ReplaceText "flask" "Flask, Git, SQL, Azure"

# --- Simple field replacements (second project block: UCLL school project) ---
ReplaceText " UCLL" " BECODE - FEENPOP"
ReplaceText "09/2019 - 2/2020" "9/2021 - 9/2021"
ReplaceText "ELIXIR/PHOENIX WEBSITE" "DATA ANONYMIZATION"

# --- Fill in the previously-empty "Tools" cell of the second project block ---
$t = $d.Tables.Item(5)
$t.Cell(15, 2).Range.Text = "SQL, Git, Streamlit, Threading, Pandas"

# --- Append 5 new project blocks (8 rows each: Environment, Methodology, Company,
#     Client, Period, Role, Tasks, Tools) right after the current last "Methodology"
#     row (row 17), i.e. before the "Management Skills" section row. ---

$blocks = @(
    @{
        Environment = ""
        Methodology = ""
        Company = "Solita"
        Client = " BECODE"
        Period = "6/2021 - 8/2021"
        Role = "BELGIUM REAL STATE PREDICTION API"
        Tasks = "The API coded in python to return the predicted price of a properties in Belgium, based on data scrapped from Immoweb from 2021. For the predictions a Linear regression was put in place to compute the relationship between several characteristics found on the sell announcement to estimate of the asking price is made. The accuracy of the model is pf 85%, which means that there is always a possibility for outliers (less then 15 %). This API has been deployed with heroku under the url: https://api-ie-predictions.herokuapp.com/ Technologies: Python, Tensorflow, PyTorch, Git, Scikit-learn, Pandas, Selenium, BeutifulSoup, HTML, Heroku"
        Tools = "Tensorflow, PyTorch, Git, Scikit-learn, Pandas"
    },
    @{
        Environment = ""
        Methodology = ""
        Company = "Solita"
        Client = " BECODE"
        Period = "5/2021 - "
        Role = "3D HOUSES VIEWER"
        Tasks = "Application coded in python to plot a house in 3D given an address. The data was collected crossing data from LIDAR satelites from Vlaanderen Overheid services and metadata, and geographical data obtained from the address through API’s servies. Technologies: Python, Git, API’s requests, Pandas, Pillow, matplotlib."
        Tools = "Git, Pandas, Pillow, matplotlib"
    },
    @{
        Environment = ""
        Methodology = ""
        Company = "Solita"
        Client = " BUSSINESS&DECISION"
        Period = "11/2021 - 3/2022"
        Role = "EDGE-COMPUTING IMAGE RECOGNITION NETWORK"
        Tasks = "Internship project. Create an internetless LAN network capable to enable communication between diferent devices exposing each one as an microservice. Two raspberries were used, first one as camera streaming service and the second as Image Recognition model host, and a phone used as endpoint to check the results of the recognition. Technologies: Python, RaspberryPi, Mimik, Yolov5, Tensorflow, MobileNetV2, RTPM"
        Tools = "RaspberryPi, Mimik, Yolov5, Tensorflow, MobileNetV2, RTPM"
    },
    @{
        Environment = ""
        Methodology = ""
        Company = "Solita"
        Client = " BECODE - FAKTION"
        Period = "10/2021 - 10/2021"
        Role = "RETAIL ANOMALY DETECTION"
        Tasks = "Application programmed in python to detect anomalies in manufacturing of dices. A Convolution Neuronal Network was trained to classify daces by its face, then a second process function computes the differences between a good manufactured one and the current one, being able to distinguish if the current dice had any anomaly or not. Technologies: Python, Git, CNN, OpenCv"
        Tools = "Git"
    },
    @{
        Environment = ""
        Methodology = ""
        Company = "Solita"
        Client = " BECODE"
        Period = "10/2021 - 10/2021"
        Role = "SIGNATURE RECOGNITION"
        Tasks = "YOLOv5 model trained in python to detect signatures on documents, it was trained with annotated documents transformed to jpg and addapting their annotations from an .xml format to a .txt normalizing and transpolating the coordinates to the yolo format. Objective: extract signatures from documents to validate legitimacy. Technologies: Python, Git, Yolov5, XML"
        Tools = "YOLOv5, Git, XML"
    }
)

$fieldOrder = @("Environment", "Methodology", "Company", "Client", "Period", "Role", "Tasks", "Tools")

# New rows must land right before the row that currently holds "Management
# Skills" (index 18). Each insertion shifts that row one index further down,
# so re-resolve the target row fresh (by its live, incrementing index) for
# every single row added -- reusing a stale row/anchor object causes the
# COM layer to keep inserting at the same spot (and thus reverses order).
$insertIndex = 18

foreach ($block in $blocks) {
    foreach ($field in $fieldOrder) {
        $targetRow = $t.Rows.Item($insertIndex)
        $newRow = $t.Rows.Add($targetRow)
        $newRow.Cells.Item(1).Range.Text = $field
        $newRow.Cells.Item(2).Range.Text = $block[$field]
        $insertIndex = $insertIndex + 1
    }
}

Write-Host "Row count after edit: " $t.Rows.Count
